$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the left table with a new "2023" column (K), mirroring the
# formatting of the existing "2022" column (J) before filling in the
# updated figures.
$ws.Range("J3:J6").Copy($ws.Range("K3:K6"))

$ws.Range("K3").Value = 2023
$ws.Range("K4").Value = 1158.5
$ws.Range("K5").Value = 559.20000000000005
$ws.Range("K6").Value = 1543
